$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Update the "panel_query_time" timestamps on the existing "data" sheet
#    (F2:F9) to the refreshed query time.
# ---------------------------------------------------------------------------
$ws1.Range("F2").Value = "2021-10-05 14:22:54.178020"
$ws1.Range("F3").Value = "2021-10-05 14:22:54.178028"
$ws1.Range("F4").Value = "2021-10-05 14:22:54.178031"
$ws1.Range("F5").Value = "2021-10-05 14:22:54.178034"
$ws1.Range("F6").Value = "2021-10-05 14:22:54.178037"
$ws1.Range("F7").Value = "2021-10-05 14:22:54.178039"
$ws1.Range("F8").Value = "2021-10-05 14:22:54.178042"
$ws1.Range("F9").Value = "2021-10-05 14:22:54.178044"

# ---------------------------------------------------------------------------
# 2) Add a new "metadata" worksheet after "data" describing the panel query
#    itself.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "metadata"

# Match the page margins used on the "data" sheet (PageSetup margins are in
# points; the "data" sheet uses 0.75in/0.75in/1in/1in/0.5in/0.5in).
$ws2.PageSetup.LeftMargin = 54
$ws2.PageSetup.RightMargin = 54
$ws2.PageSetup.TopMargin = 72
$ws2.PageSetup.BottomMargin = 72
$ws2.PageSetup.HeaderMargin = 36
$ws2.PageSetup.FooterMargin = 36

# Copy the bold/bordered/centered header formatting used on the "data" sheet
# onto the new header row and the index column, re-using the existing style
# rather than inventing a new one.
$ws1.Range("B1").Copy()
$ws2.Range("B1:G1").PasteSpecial(-4122)
$ws2.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws2.Range("B1").Value = "data_name"
$ws2.Range("C1").Value = "data_id"
$ws2.Range("D1").Value = "data_version"
$ws2.Range("E1").Value = "data_version_created"
$ws2.Range("F1").Value = "panel_query_time"
$ws2.Range("G1").Value = "panel_get_request"

$ws2.Range("A2").Value = 0
$ws2.Range("B2").Value = "Surfactant deficiency"
$ws2.Range("C2").Value = 551

# "1.9" must be stored as text, not as the number 1.9 - round-trip it through
# a formula then paste-special as a value so it lands as a plain string with
# no extra number-format / style baggage.
$ws2.Range("D2").Formula = '="1.9"'
$ws2.Range("D2").Copy()
$ws2.Range("D2").PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws2.Range("E2").Value = "2021-03-17T14:01:07.529142Z"
$ws2.Range("F2").Value = "2021-10-05 14:22:54.174347"
$ws2.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/551/?format=json"

# Keep "data" as the active sheet/tab, matching the original activeTab.
$ws1.Activate()
